$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reorder/rename columns ---
# lang_code | code | name | is_active
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "is_active"

# --- Data rows 2-7 ---
$data = @(
    @("eng", "MLE", "Male",    $true),
    @("eng", "FLE", "Female",  $true),
    @("eng", "OTH", "Others",  $false),
    @("fra", "MLE", "Mâle",    $true),
    @("fra", "FLE", "Femelle", $true),
    @("fra", "OTH", "Dautres", $false)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
    $ws.Range("D$row").Value = $vals[3]
}

# Column A (rows 2-7) takes on the same style as the header cells
# (bold, bordered, centered) -- copy that formatting from A1.
$ws.Range("A1").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A1").Select()
